$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the numeric-looking text cells in columns D and G keep their original
# text (string) representation -- force text number format before assigning
# so Excel does not reinterpret these values as actual numbers (which would
# drop significant trailing zeros, e.g. "5.360" -> 5.36).
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "246.03"
$ws.Range("G2").Value = "18"
$ws.Range("D3").Value = "22.07"
$ws.Range("G3").Value = "18"
$ws.Range("D4").Value = "5.360"
$ws.Range("G4").Value = "18"
$ws.Range("D5").Value = "0.05854"
$ws.Range("G5").Value = "18"
$ws.Range("G6").Value = "18"
$ws.Range("D7").Value = "6.373"
$ws.Range("G7").Value = "18"
$ws.Range("D8").Value = "0.8140"
$ws.Range("G8").Value = "18"
$ws.Range("D9").Value = "1.013"
$ws.Range("G9").Value = "18"
$ws.Range("D10").Value = "0.1420"
$ws.Range("G10").Value = "18"
$ws.Range("D11").Value = "0.04156"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G11").Value = "18"
$ws.Range("D12").Value = "0.07369"
$ws.Range("G12").Value = "18"
$ws.Range("D13").Value = "0.03033"
$ws.Range("G13").Value = "18"
$ws.Range("D14").Value = "4.174"
$ws.Range("G14").Value = "18"
$ws.Range("D15").Value = "0.09398"
$ws.Range("G15").Value = "18"
$ws.Range("D16").Value = "0.001593"
$ws.Range("G16").Value = "18"
$ws.Range("D17").Value = "0.04803"
$ws.Range("G17").Value = "18"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "0.006065"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("G18").Value = "18"
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D19").Value = "0.004080"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("G19").Value = "18"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Value = "0.0009831"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("G20").Value = "18"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").Value = "0.0001500"
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("G21").Value = "18"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").Value = "3.707"
$ws.Range("E22").Value = "21LEOLEO"
$ws.Range("G22").Value = "18"
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D23").Value = "2.232"
$ws.Range("E23").Value = "22BTSETokenBTSE"
$ws.Range("G23").Value = "18"
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D24").Value = "0.0005890"
$ws.Range("E24").Value = "23OneONE"
$ws.Range("G24").Value = "18"
$ws.Range("D25").Value = "0.3239"
$ws.Range("G25").Value = "18"
$ws.Range("G26").Value = "18"
$ws.Range("G27").Value = "18"
$ws.Range("G28").Value = "18"
$ws.Range("G29").Value = "18"
$ws.Range("G30").Value = "18"
$ws.Range("G31").Value = "18"
$ws.Range("G32").Value = "18"
$ws.Range("G33").Value = "18"
$ws.Range("G34").Value = "18"
$ws.Range("G35").Value = "18"
$ws.Range("G36").Value = "18"
$ws.Range("G37").Value = "18"
$ws.Range("G38").Value = "18"
$ws.Range("G39").Value = "18"
$ws.Range("D40").Value = "0.03861"
$ws.Range("G40").Value = "18"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "0.006369"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("G41").Value = "18"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "0.1074"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("G42").Value = "18"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "0.002600"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("G43").Value = "18"
$ws.Range("G44").Value = "18"
$ws.Range("D45").Value = "0.00005627"
$ws.Range("G45").Value = "18"
$ws.Range("G46").Value = "18"
$ws.Range("D47").Value = "0.9399"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINBestin24h"
$ws.Range("G47").Value = "18"
$ws.Range("D48").Value = "0.07344"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
$ws.Range("G48").Value = "18"
$ws.Range("G49").Value = "18"
$ws.Range("G50").Value = "18"
$ws.Range("G51").Value = "18"
